# 自动更新Excel文件 - 2025-10-27 23:12:01
# Daily countdown update: for every data row (2..99), the "剩余" (remaining
# days, column E) counts down by 1 day. When a row's remaining count has
# reached 1 (its last day), the cycle restarts: "剩余" resets to the row's
# "总天" (total days, column D) and "开始时间" (start date, column F) is
# bumped to the new cycle's start date (20251028).
# Row 36 is skipped: its "开始时间" value is corrupted (202510929, not a
# valid yyyymmdd date), so the source data for that row was left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newStartDate = 20251028

for ($row = 2; $row -le 99; $row++) {
    if ($row -eq 36) {
        continue
    }

    $totalDays = $ws.Cells.Item($row, 4).Value2
    $remaining = $ws.Cells.Item($row, 5).Value2

    if ($remaining -eq 1) {
        $ws.Cells.Item($row, 5).Value = $totalDays
        $ws.Cells.Item($row, 6).Value = $newStartDate
    } else {
        $ws.Cells.Item($row, 5).Value = $remaining - 1
    }
}
